# Nieuwe werkende sjablonen + temp cache fix
#
# This applies the "nieuwe werkende sjablonen" edit: the single title
# slide gets re-pointed at the (already present) "Titeldia" layout
# placeholders, and four brand-new slides are appended, one per
# remaining custom layout, each carrying just the bare photo
# placeholders that that layout defines.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layouts = $master.CustomLayouts

# --- Slide 1: swap the generic ctrTitle/subTitle placeholders for the
#     dedicated "Titeldia" layout's own placeholders -------------------
$s1 = $p.Slides.Item(1)
$titleLayout = $layouts.Item(1)            # "Titeldia"
$s1.CustomLayout = $titleLayout

# CustomLayout assignment adds the layout's placeholders alongside the
# two legacy ones (their ph type/idx no longer matches, so they are not
# reused) -- drop the two stale shapes, leaving only the freshly added
# "Naam+geboortedatum+overleidensdatum" / "foto_01" placeholders.
$s1.Shapes.Item(1).Delete()
$s1.Shapes.Item(1).Delete()
$s1.Shapes.Item(2).Name = "foto01"

# --- Slides 2-5: one new slide per remaining layout, each with just its
#     picture placeholders (no text typed into any of them) -----------
$s2 = $p.Slides.AddSlide(2, $layouts.Item(2))   # "Titel en object"
$s2.Shapes.Item(1).Name = "foto03"
$s2.Shapes.Item(2).Name = "foto02"

$s3 = $p.Slides.AddSlide(3, $layouts.Item(3))   # "Sectiekop"
$s3.Shapes.Item(1).Name = "foto06"
$s3.Shapes.Item(2).Name = "foto05"
$s3.Shapes.Item(3).Name = "foto04"

$s4 = $p.Slides.AddSlide(4, $layouts.Item(4))   # "Inhoud van twee"
$s4.Shapes.Item(1).Name = "foto08"
$s4.Shapes.Item(2).Name = "foto07"

$s5 = $p.Slides.AddSlide(5, $layouts.Item(5))   # "Vergelijking"
$s5.Shapes.Item(1).Name = "foto09"

Write-Host "Slides: $($p.Slides.Count)"
